# DataSources.xlsx - "Unrestricted datapipeline updated - need source change or update"
#
# 1) Simplify the SRC column (C) values from "<name>_actual"/"<name>_forecast"
#    down to just the series name (FERT / FPI / GDP / NG / USDEUR).
# 2) Add a new "INTERVAL" column (I) to Table2 carrying the Q/M/A interval
#    code that used to be implied by the _actual/_forecast suffix.
# 3) Widen column C a bit so the new, slightly longer values still read well.
# 4) Leave the view zoomed out a touch and the selection parked on I4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Collapse the SRC column values -------------------------------------
$ws.Range("C2").Value = "FERT"
$ws.Range("C3").Value = "FERT"
$ws.Range("C4").Value = "FPI"
$ws.Range("C5").Value = "FPI"
$ws.Range("C6").Value = "GDP"
$ws.Range("C7").Value = "GDP"
$ws.Range("C8").Value = "NG"
$ws.Range("C9").Value = "NG"
$ws.Range("C10").Value = "USDEUR"
$ws.Range("C11").Value = "USDEUR"

# --- 2) Add the INTERVAL column to the table --------------------------------
$lo = $ws.ListObjects.Item(1)
$newCol = $lo.ListColumns.Add()

$ws.Range("I1").Value = "INTERVAL"
# give the new header the same plain bordered look as the rest of row 1
# (thin border on left/right/bottom, no top - matches the other header cells)
$hdr = $ws.Range("I1")
$hdr.Borders.LineStyle = 1
$hdr.Borders.Item(8).LineStyle = -4142

$ws.Range("I2").Value = "A"
$ws.Range("I3").Value = "A"
$ws.Range("I4").Value = "M"
# I5 intentionally left blank (matches source data)
$ws.Range("I6").Value = "Q"
$ws.Range("I7").Value = "Q"
$ws.Range("I8").Value = "M"
$ws.Range("I9").Value = "M"
$ws.Range("I10").Value = "M"
$ws.Range("I11").Value = "A"

# --- 3) Widen column C --------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 13.6

# --- 4) View tweaks ----------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 85
$ws.Range("I4").Select()
